$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -6.782899999999996
$ws.Range("A4").Value = -20.83639999999998
$ws.Range("A6").Value = -23.06720000000001
$ws.Range("A7").Value = -21.90619999999998
$ws.Range("B7").Value = 5.303200000000003
$ws.Range("A8").Value = -22.25260000000001
$ws.Range("B11").Value = 5.647699999999998
$ws.Range("B12").Value = 5.962599999999997
$ws.Range("C12").Value = -11.28900000000001
$ws.Range("D12").Value = -8.043500000000003
$ws.Range("C13").Value = -12.85199999999999
$ws.Range("D13").Value = -8.815000000000005
$ws.Range("C14").Value = -14.649
$ws.Range("B15").Value = 5.303999999999994
$ws.Range("A16").Value = -21.48009999999998
$ws.Range("C16").Value = -11.5799
$ws.Range("C19").Value = -12.95300000000001
$ws.Range("A20").Value = -22.5585
$ws.Range("B20").Value = 5.000999999999999
$ws.Range("C20").Value = -13.52589999999999
$ws.Range("A21").Value = -22.2801
$ws.Range("B21").Value = 5.361699999999997
$ws.Range("B22").Value = 9.0807
$ws.Range("C22").Value = -12.4142
$ws.Range("D22").Value = -8.178499999999998
$ws.Range("B23").Value = 9.426799999999997
$ws.Range("D25").Value = -8.3314
$ws.Range("A28").Value = -22.08519999999999
$ws.Range("A29").Value = -21.63710000000001
$ws.Range("B29").Value = 5.074199999999999
$ws.Range("D29").Value = -6.949199999999995
$ws.Range("A30").Value = -21.61060000000001
$ws.Range("A32").Value = -21.273
$ws.Range("B34").Value = 8.563200000000002
$ws.Range("D34").Value = -7.972700000000003
$ws.Range("C36").Value = -13.0074
$ws.Range("A40").Value = -18.9471
$ws.Range("B42").Value = 10.5169
$ws.Range("B43").Value = 4.523799999999992
$ws.Range("C43").Value = -14.42969999999998
$ws.Range("D43").Value = -8.647899999999995
$ws.Range("B44").Value = 4.726999999999999
$ws.Range("B45").Value = 5.194100000000001
$ws.Range("A46").Value = -21.55830000000001
$ws.Range("B46").Value = 5.634000000000002
$ws.Range("C46").Value = -11.28400000000001
$ws.Range("D48").Value = -8.131300000000003
$ws.Range("B50").Value = 4.847799999999998
$ws.Range("C50").Value = -13.64679999999999
$ws.Range("A51").Value = -22.27029999999999
$ws.Range("B51").Value = 5.328799999999998
$ws.Range("A52").Value = -22.2247
$ws.Range("A57").Value = -23.0029
$ws.Range("B57").Value = 4.833599999999997
$ws.Range("A59").Value = -22.42179999999998
$ws.Range("D60").Value = -8.626999999999999
$ws.Range("A62").Value = -22.31499999999999
$ws.Range("B65").Value = 5.819300000000001
$ws.Range("A66").Value = -21.50269999999999
$ws.Range("B66").Value = 4.876299999999995
$ws.Range("B67").Value = 5.303599999999997
$ws.Range("D68").Value = -6.689299999999998
$ws.Range("D70").Value = -6.774799999999996
$ws.Range("D71").Value = -7.416099999999996
$ws.Range("A73").Value = -19.12009999999998
$ws.Range("D73").Value = -8.212400000000001
$ws.Range("A74").Value = -21.88209999999999
$ws.Range("C76").Value = -12.028
$ws.Range("A77").Value = -20.28169999999999
$ws.Range("D78").Value = -8.296000000000003
$ws.Range("B79").Value = 9.964600000000011
$ws.Range("B84").Value = 5.0854
$ws.Range("B87").Value = 4.682199999999996
$ws.Range("D87").Value = -8.62539999999999
$ws.Range("A92").Value = -21.57400000000002
$ws.Range("B92").Value = 4.832099999999998
$ws.Range("D92").Value = -6.279199999999999
$ws.Range("C95").Value = -11.08660000000001
$ws.Range("B97").Value = 6.431799999999996
$ws.Range("C97").Value = -11.27660000000001
$ws.Range("C99").Value = -12.04080000000001
$ws.Range("A100").Value = -22.4204
$ws.Range("D101").Value = -7.980300000000003
